$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.855.94'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '1.903.36'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5023'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3815'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07287'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9109'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.51%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07650'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = '1.902.28'
$ws.Range("E13").Value = '  -0.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.484'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.613'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008718'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = '27.880.64'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.52'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.151'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '154.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.862'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.229'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '115.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.933'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08971'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.212'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.53%  '
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7707'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.643'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02056'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.50%  '
$ws.Range("E36").Value = '  -1.93%  '
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5540'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.014'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05275'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.986'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.550'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1524'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '111.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.61'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4802'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.002'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '67.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06084'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9007'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.46%  '
